$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Reorder: "Plasmid" moves before "PlasmidBatch"
#    (SOP, PlasmidBatch, Plasmid, ...) -> (SOP, Plasmid, PlasmidBatch, ...)
# ---------------------------------------------------------------------------
$plasmid = $wb.Worksheets.Item("Plasmid")
$plasmidBatch = $wb.Worksheets.Item("PlasmidBatch")
$plasmid.Move($plasmidBatch)

# NOTE: after Move() re-orders the tab strip, previously-held worksheet
# variables track tab *position*, not identity, so re-fetch by name before
# using them again.
$plasmid = $wb.Worksheets.Item("Plasmid")

# ---------------------------------------------------------------------------
# 2. "Plasmid" header row: drop "Location", keep "Batches" as last column
#    name, responsible, CommonName, Usage, Antibiotic, Location, Batches
#    -> name, responsible, CommonName, Usage, Antibiotic, Batches
# ---------------------------------------------------------------------------
$plasmid.Cells.Item(1, 6).Value = "Batches"
$plasmid.Cells.Item(1, 7).ClearContents()

# ---------------------------------------------------------------------------
# 3. "Strain" header row: drop TubesLeft / TubeVolume / Location, add "Batches"
#    ... CatalogNo, TubesLeft, TubeVolume, Location (A1:O1)
#    -> ... CatalogNo, Batches (A1:M1)
# ---------------------------------------------------------------------------
$strain = $wb.Worksheets.Item("Strain")
$strain.Cells.Item(1, 13).Value = "Batches"
$strain.Cells.Item(1, 14).ClearContents()
$strain.Cells.Item(1, 15).ClearContents()

# ---------------------------------------------------------------------------
# 4. Rename "CellLine" -> "StrainBatch" and replace its field columns
#    name, responsible, CommonName, Usage, Species, Genotype, Location
#    -> name, responsible, TubesLeft, Barcode, Location, SequenceVerified, Parent
# ---------------------------------------------------------------------------
$strainBatch = $wb.Worksheets.Item("CellLine")
$strainBatch.Name = "StrainBatch"
$strainBatch.Cells.Item(1, 3).Value = "TubesLeft"
$strainBatch.Cells.Item(1, 4).Value = "Barcode"
$strainBatch.Cells.Item(1, 5).Value = "Location"
$strainBatch.Cells.Item(1, 6).Value = "SequenceVerified"
$strainBatch.Cells.Item(1, 7).Value = "Parent"

# ---------------------------------------------------------------------------
# 5. New sheet "CellLine" right after "StrainBatch"
#    name, responsible, CommonName, Usage, Species, Genotype, Batches
# ---------------------------------------------------------------------------
$cellLine = $wb.Worksheets.Add($null, $strainBatch)
$cellLine.Name = "CellLine"
$cellLineHeaders = @("name", "responsible", "CommonName", "Usage", "Species", "Genotype", "Batches")
for ($i = 0; $i -lt $cellLineHeaders.Length; $i++) {
    $cellLine.Cells.Item(1, $i + 1).Value = $cellLineHeaders[$i]
}

# ---------------------------------------------------------------------------
# 6. New sheet "CellLineBatch" right after "CellLine"
#    name, responsible, Location, Barcode, TubesLeft, Mycoplasma, Parent
# ---------------------------------------------------------------------------
$cellLineBatch = $wb.Worksheets.Add($null, $cellLine)
$cellLineBatch.Name = "CellLineBatch"
$cellLineBatchHeaders = @("name", "responsible", "Location", "Barcode", "TubesLeft", "Mycoplasma", "Parent")
for ($i = 0; $i -lt $cellLineBatchHeaders.Length; $i++) {
    $cellLineBatch.Cells.Item(1, $i + 1).Value = $cellLineBatchHeaders[$i]
}

# ---------------------------------------------------------------------------
# 7. New sheet "CultureMedia" right after "CellLineBatch"
#    name, responsible, ProductName, Vendor, CatalogNo, Batches
# ---------------------------------------------------------------------------
$cultureMedia = $wb.Worksheets.Add($null, $cellLineBatch)
$cultureMedia.Name = "CultureMedia"
$cultureMediaHeaders = @("name", "responsible", "ProductName", "Vendor", "CatalogNo", "Batches")
for ($i = 0; $i -lt $cultureMediaHeaders.Length; $i++) {
    $cultureMedia.Cells.Item(1, $i + 1).Value = $cultureMediaHeaders[$i]
}

# ---------------------------------------------------------------------------
# 8. New sheet "CultureMediaBatch" right after "CultureMedia"
#    name, responsible, Location, Barcode, TubesLeft, Mycoplasma, Parent
# ---------------------------------------------------------------------------
$cultureMediaBatch = $wb.Worksheets.Add($null, $cultureMedia)
$cultureMediaBatch.Name = "CultureMediaBatch"
$cultureMediaBatchHeaders = @("name", "responsible", "Location", "Barcode", "TubesLeft", "Mycoplasma", "Parent")
for ($i = 0; $i -lt $cultureMediaBatchHeaders.Length; $i++) {
    $cultureMediaBatch.Cells.Item(1, $i + 1).Value = $cultureMediaBatchHeaders[$i]
}

# Fermentation sheet stays last, content unchanged.
